$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Correct" (Right) marking value for the Total row's Marking row (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update the Total correct marks (B12: 72 -> 120)
$ws.Range("B12").Value = 120

# Update the "Corr/total marks" summary text (E12: "69/84" -> "120/140")
$ws.Range("E12").Value = "120/140"
